# Generate Report for Handback
#
# This mirrors the "handback" report-generation step of the OpenLocalization
# CI tool: once a translated xliff has round-tripped back and is in sync
# with the en-US source, the status columns move from "Ready for handoff"
# to "Handed back: in sync with en-US", and the per-language worksheets get
# their "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns populated (plus a hyperlink on the target-file cell,
# matching the existing source-file hyperlink). The wider text that now
# lives in those cells means the report also widens a few columns so the
# values aren't clipped.

$wb = $excel.ActiveWorkbook
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both language status columns flip to "handed back".
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

# Widen the (now longer) status columns so the text isn't truncated.
$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------
# Helper data shared by both language sheets: source-file names/urls.
# ---------------------------------------------------------------------
$url091878 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/895dedfc25851b551e2e65d7f6ee8d4a7e22e44b/e2e/091878de-e876-46e4-bd89-d8ae4ae75a69.md"
$url2b57fc = "https://github.com/OpenLocalizationTestOrg/oltest/blob/895dedfc25851b551e2e65d7f6ee8d4a7e22e44b/e2e/2b57fcd6-56ee-4faa-971a-ba26295e6feb.md"
$name091878 = "091878de-e876-46e4-bd89-d8ae4ae75a69.md"
$name2b57fc = "2b57fcd6-56ee-4faa-971a-ba26295e6feb.md"

# ---------------------------------------------------------------------
# zh-cn sheet: fill in Latest Target File (I), Latest Handback File (J)
# and Latest Handback DateTime (K) for both rows.
# ---------------------------------------------------------------------
$zhcn.Range("J2").Value = "091878de-e876-46e4-bd89-d8ae4ae75a69.1986e1b365d920c66086e70db56edeb7c2c6fc86.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-12 12:31:34"

$zhcn.Range("J3").Value = "2b57fcd6-56ee-4faa-971a-ba26295e6feb.da20f53a1a52a2cb3d59594bd4f716a9c20d0fba.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-12 12:31:34"

# Latest Target File (I) mirrors the source-file hyperlink/name.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $url091878, "", "", $name091878)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $url091878, "", "", $name091878)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $url2b57fc, "", "", $name2b57fc)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $url2b57fc, "", "", $name2b57fc)

$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: same shape, different language-specific file names and
# its own handback timestamp.
# ---------------------------------------------------------------------
$dede.Range("J2").Value = "091878de-e876-46e4-bd89-d8ae4ae75a69.1986e1b365d920c66086e70db56edeb7c2c6fc86.de-de.xlf"
$dede.Range("K2").Value = "2016-08-12 12:31:43"

$dede.Range("J3").Value = "2b57fcd6-56ee-4faa-971a-ba26295e6feb.da20f53a1a52a2cb3d59594bd4f716a9c20d0fba.de-de.xlf"
$dede.Range("K3").Value = "2016-08-12 12:31:43"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $url091878, "", "", $name091878)
$dede.Hyperlinks.Add($dede.Range("I2"), $url091878, "", "", $name091878)
$dede.Hyperlinks.Add($dede.Range("A3"), $url2b57fc, "", "", $name2b57fc)
$dede.Hyperlinks.Add($dede.Range("I3"), $url2b57fc, "", "", $name2b57fc)

$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17
